# Updates cryptos list (price/volume refresh + WrappedEther/Polkadot row swap),
# matching the Wed Aug 2 04:30:36 UTC 2023 GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text such as
# "29.631.06", "1.000" or "160.00". Excel's COM Value setter would silently
# coerce these to numbers/dates and strip the formatting, so each such cell
# is forced to Text format ("@") immediately before its value is written.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.631.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.862.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7008'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07748'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3074'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.69'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.873.87'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.167'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.49'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6940'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.603.96'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008375'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.111.50'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '242.53'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.79'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.637'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1513'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.931'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.00'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.542'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.269'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.187'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.194'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05108'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7862'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.906'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.158'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.333.16'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +10.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01881'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.14%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9576'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +14.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.71'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000126'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.798'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.009.46'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5220'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.23'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.790'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.020'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.56%  '
